# InputExport.xlsx — refresh the export with a new single-line snapshot.
# The sheet previously listed 15 part rows (rows 2-16); the new export
# only contains one data row, so rows 3-16 are removed and row 2 is
# overwritten with the latest record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old extra data rows (3-16), leaving just the header + 1 row.
$ws.Rows("3:16").Delete()

# Overwrite row 2 with the new record.
$ws.Range("A2").Value = "T109.410.11.072.00"

$ws.Range("B2").Value = 1
$ws.Range("B2").Font.Size = 11
$ws.Range("B2").Font.Name = "新細明體"

$ws.Range("C2").Value = "23/11/2023 POP"

$ws.Range("D2").Value = 1119

# The averaging helper formula's reference became stale during the
# refresh, so it now evaluates to a #REF! error.
$ws.Range("E2").Formula = "=#REF!/B2"

# Leave the cursor where the next entry would be typed.
$ws.Range("B3").Select()
